# Vendors.xlsx edit: switch "current invoice" row from CNA Insurance (row 6)
# to Blue Cross Blue Shield of Illinois (row 3) and ComEd (row 8), bump the
# UPS invoice number, add a "source file" marker column, and nudge the
# tracking counter in I1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the running counter in I1 and tag the header row with the script
# that produced this pass (new column J).
$ws.Range("I1").Value = 33335
$ws.Range("J1").Value = "bill_entry_and_payment.py"

# Blue Cross Blue Shield of Illinois (row 3): mark for payment, set bill amount.
$ws.Range("B3").Value = "x"
$ws.Range("D3").Value = 4670.8999999999996

# CNA Insurance (row 6): no longer flagged / billed this pass.
$ws.Range("B6").ClearContents()
$ws.Range("D6").ClearContents()

# ComEd (row 8): mark for payment, set bill amount.
$ws.Range("B8").Value = "x"
$ws.Range("D8").Value = 409.98

# UPS (row 22): updated invoice number.
$ws.Range("C22").Value = "00007RW518503"

# Leave the selection where the last edit happened.
$ws.Range("D8").Select()
